# "Updating all test cases" - append one more test-case row (Index 20,
# Reference_Key "ABC_20", Status "sample_status_20") to the bottom of the
# existing table on Sheet1 (which currently runs from row 2 through row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 21

$ws.Cells.Item($newRow, 1).Value = 20
$ws.Cells.Item($newRow, 2).Value = "ABC_20"
$ws.Cells.Item($newRow, 3).Value = "sample_status_20"
